$wb = $excel.ActiveWorkbook

# --- Sheet1: move selection from A14 to A13 (tabSelected will move off Sheet1) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A13").Select()

# --- Sheet2: becomes the active tab, gains test-case formulas in A1:A6 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Formula = "=1+1"
$ws2.Range("A2").Formula = "=1-1"
$ws2.Range("A3").Formula = "=1+1"
$ws2.Range("A4").Formula = "=1--1"
$ws2.Range("A5").Formula = "=1+-1"
$ws2.Range("A6").Formula = "=1-1"

# Activating Sheet2 last makes it the selected/active tab (workbook activeTab=1,
# sheet2 tabSelected=1, sheet1 loses tabSelected) while keeping the A13 selection on Sheet1.
$ws2.Activate()
